# Revert to previous version of Master_Fuel_Sector_List.xlsx.
# - Fixes activity/units for a few sectors that had been miscoded as
#   GDP/B2005USD instead of Energy_Combustion/kt.
# - Restores three missing sector rows (11A/11B/11C - natural sources).
# - Removes the unused Hyperlink / Followed Hyperlink cell styles.
# - Restores prior view state (Fuels tab active, Sectors selection at D24).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sectors")
$ws2 = $wb.Worksheets.Item("Fuels")

# --- Data fixes on "Sectors": these rows should use Energy_Combustion/kt,
#     not GDP/B2005USD ---
$ws1.Range("B5").Value  = "Energy_Combustion"
$ws1.Range("C5").Value  = "kt"

$ws1.Range("B30").Value = "Energy_Combustion"
$ws1.Range("C30").Value = "kt"

$ws1.Range("B31").Value = "Energy_Combustion"
$ws1.Range("C31").Value = "kt"

$ws1.Range("B32").Value = "Energy_Combustion"
$ws1.Range("C32").Value = "kt"

# --- Restore the three missing natural-source sector rows ---
$ws1.Range("A57").Value = "11A_Volcanoes"
$ws1.Range("B57").Value = "GDP"
$ws1.Range("C57").Value = "B2005USD"

$ws1.Range("A58").Value = "11B_Forest-fires"
$ws1.Range("B58").Value = "GDP"
$ws1.Range("C58").Value = "B2005USD"

$ws1.Range("A59").Value = "11C_Other-natural"
$ws1.Range("B59").Value = "GDP"
$ws1.Range("C59").Value = "B2005USD"

# --- Remove unused legacy hyperlink cell styles (two duplicate pairs) ---
$wb.Styles.Item("Followed Hyperlink").Delete()
$wb.Styles.Item("Followed Hyperlink").Delete()
$wb.Styles.Item("Hyperlink").Delete()
$wb.Styles.Item("Hyperlink").Delete()

# --- Restore prior view/selection state ---
$ws1.Range("D24").Select()
$ws2.Range("C30").Select()

# "Fuels" is the active tab in the restored version.
$ws2.Activate()
